$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.240.99"
$ws.Range("E2").Value = "  -1.95%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.501.89"
$ws.Range("E3").Value = "  -4.88%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.76"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.59"
$ws.Range("E6").Value = "  +1.16%  "
$ws.Range("E8").Value = "  -2.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.501.92"
$ws.Range("E9").Value = "  -4.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.137"
$ws.Range("E10").Value = "  -2.05%  "
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("E12").Value = "  -4.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.09"
$ws.Range("E13").Value = "  -2.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.41"
$ws.Range("E14").Value = "  -4.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.939.56"
$ws.Range("E15").Value = "  -5.50%  "
$ws.Range("E16").Value = "  -4.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.009.29"
$ws.Range("E17").Value = "  -2.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.497.59"
$ws.Range("E18").Value = "  -5.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.20"
$ws.Range("E19").Value = "  -6.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.64"
$ws.Range("E20").Value = "  -5.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "344.93"
$ws.Range("E21").Value = "  -3.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.18"
$ws.Range("E22").Value = "  -3.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.55"
$ws.Range("E23").Value = "  -2.72%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("E25").Value = "  -0.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "69.11"
$ws.Range("E26").Value = "  -0.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.82"
$ws.Range("E27").Value = "  -5.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.623.55"
$ws.Range("E29").Value = "  -5.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0963"
$ws.Range("E30").Value = "  -4.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "525.24"
$ws.Range("E31").Value = "  -4.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.06"
$ws.Range("E32").Value = "  +1.56%  "
$ws.Range("E33").Value = "  -3.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.82"
$ws.Range("E34").Value = "  -4.15%  "
$ws.Range("E35").Value = "  -3.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "157.71"
$ws.Range("E37").Value = "  +0.69%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.45"
$ws.Range("E38").Value = "  -4.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.49"
$ws.Range("E39").Value = "  -2.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.33"
$ws.Range("E40").Value = "  +0.23%  "
$ws.Range("E41").Value = "  -3.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.77"
$ws.Range("E42").Value = "  -2.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.05"
$ws.Range("E43").Value = "  -3.52%  "
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.43"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "146.42"
$ws.Range("E46").Value = "  -4.29%  "
$ws.Range("E47").Value = "  -4.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.66"
$ws.Range("E48").Value = "  -3.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.72"
$ws.Range("E49").Value = "  +1.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0269"
$ws.Range("E50").Value = "  -9.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0751"
$ws.Range("E51").Value = "  -2.63%  "
